$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns are stored as literal text in the
# source data (inline strings), including values that look numeric (e.g.
# "319.18", "0.4385"). Plain Value assignment would make Excel coerce those
# into real numbers, so we force the Text number format first for any D-cell
# whose new value is a bare number (it is a no-op, COM-faithful step for the
# multi-dot "28.228.84"-style values, which never parse as numbers anyway).

$ws.Range("D2").Value = "28.228.84"
$ws.Range("E2").Value = "  -2.39%  "

$ws.Range("D3").Value = "1.863.54"
$ws.Range("E3").Value = "  -1.89%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.18"
$ws.Range("E5").Value = "  -1.58%  "

$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4385"
$ws.Range("E7").Value = "  -4.46%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3721"
$ws.Range("E8").Value = "  -2.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07540"
$ws.Range("E9").Value = "  -2.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9393"
$ws.Range("E10").Value = "  -3.51%  "

$ws.Range("D12").Value = "1.884.44"
$ws.Range("E12").Value = "  -1.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.731"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.448"
$ws.Range("E14").Value = "  -3.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06866"
$ws.Range("E15").Value = "  -2.60%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "82.20"
$ws.Range("E17").Value = "  -1.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009118"
$ws.Range("E18").Value = "  -3.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  -0.18%  "

$ws.Range("E20").Value = "  -3.61%  "

$ws.Range("D21").Value = "28.221.47"
$ws.Range("E21").Value = "  -2.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.165"
$ws.Range("E22").Value = "  -2.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.77"
$ws.Range("E23").Value = "  -0.77%  "

$ws.Range("D24").Value = "2.093.76"
$ws.Range("E24").Value = "  -2.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.039"
$ws.Range("E25").Value = "  -2.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.84"
$ws.Range("E26").Value = "  -2.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.46"
$ws.Range("E27").Value = "  -2.96%  "

$ws.Range("E28").Value = "  -4.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.31"
$ws.Range("E29").Value = "  -2.57%  "

$ws.Range("E30").Value = "  -5.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09061"
$ws.Range("E31").Value = "  -1.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8031"
$ws.Range("E32").Value = "  -6.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.859"
$ws.Range("E33").Value = "  -4.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.171"
$ws.Range("E34").Value = "  -5.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.955"
$ws.Range("E35").Value = "  -1.47%  "

$ws.Range("E36").Value = "  -0.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.123"
$ws.Range("E37").Value = "  -1.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05467"
$ws.Range("E38").Value = "  -3.46%  "

$ws.Range("E39").Value = "  -3.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.975"
$ws.Range("E40").Value = "  +7.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.145"
$ws.Range("E41").Value = "  -3.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5260"
$ws.Range("E42").Value = "  -3.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1674"
$ws.Range("E43").Value = "  -4.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.782"
$ws.Range("E44").Value = "  -5.13%  "

$ws.Range("E45").Value = "  +1.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06771"
$ws.Range("E46").Value = "  -0.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4873"
$ws.Range("E47").Value = "  -5.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000002516"
$ws.Range("E48").Value = "  -2.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.57"
$ws.Range("E49").Value = "  -5.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "107.71"
$ws.Range("E50").Value = "  -2.11%  "

$ws.Range("E51").Value = "  -4.63%  "
